$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update metrics for the last row (year 2025) with refreshed data
$ws.Range("C8").Value = 1456
$ws.Range("D8").Value = 227
$ws.Range("E8").Value = 1229
$ws.Range("F8").Value = 9.310910582444627
$ws.Range("G8").Value = 84.40934065934066
$ws.Range("H8").Value = 15.59065934065934
